$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# 1) "N Felder zurück" -> "N Drinks oder n Felder zurück"
$d.Content.Find.Execute("N Felder zurück", $true, $false, $false, $false, $false,
                         $true, 1, $false, "N Drinks oder n Felder zurück", 2)

# 2) A new list item "Armdrückduel" is added after "Würfeln+trinken". In the
#    original document the trailing "_GoBack" bookmark (Word's "last edit
#    location" marker) sits inside the "Würfeln+trinken" paragraph; after the
#    edit it has moved into the new "Armdrückduel" paragraph, so we rebuild
#    both paragraphs to match.

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Keep the paragraph's own revision-save-id attributes (w:rsidR/.../w:rsidP)
# intact instead of letting them get reset by the rewrite below.
$attrs = ""
$oxml = $lastRange.WordOpenXML
if ($oxml -match '<w:p (?:w14:paraId="[^"]*" w14:textId="[^"]*" )?([^>]*)><w:pPr><w:pStyle w:val="Listenabsatz"') {
    $attrs = $matches[1]
}

# Rewrite "Würfeln+trinken" paragraph without the bookmark.
$wuerfelnXml = '<w:p ' + $attrs + ' xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Würfeln+trinken</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$lastRange.InsertXML($wuerfelnXml)

# Append the new "Armdrückduel" list paragraph, now carrying the bookmark.
$endRange = $d.Content
$endRange.Collapse(0)
$armdrueckXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Armdrückduel</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p>'
$endRange.InsertXML($armdrueckXml)
